$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for provinces that no longer exist (old rows 63-69),
# shrinking the table from 68 provinces down to 61.
$ws.Range("A63:A69").EntireRow.Delete()

# Rewrite the province list (column A) and the three hazard-tier columns
# (B = risk_to_assets, C = resilience, D = risk) for every remaining row,
# now alphabetically sorted and including the updated set of provinces
# (Batanes, Misamis Occidental, Negros Oriental, Siquijor added; the
# Del Norte/Del Sur split provinces, North Cotabato, Samar and Surigao
# Del Norte removed), populated with illustrative coastal flood data.
$ws.Cells.Item(2, 1).Value = "Abra"
$ws.Cells.Item(2, 2).Value = "High"
$ws.Cells.Item(2, 3).Value = "Mid"
$ws.Cells.Item(2, 4).Value = "High"
$ws.Cells.Item(3, 1).Value = "Aklan"
$ws.Cells.Item(3, 2).Value = "High"
$ws.Cells.Item(3, 3).Value = "Mid"
$ws.Cells.Item(3, 4).Value = "High"
$ws.Cells.Item(4, 1).Value = "Albay"
$ws.Cells.Item(4, 2).Value = "Low"
$ws.Cells.Item(4, 3).Value = "Low"
$ws.Cells.Item(4, 4).Value = "Low"
$ws.Cells.Item(5, 1).Value = "Antique"
$ws.Cells.Item(5, 2).Value = "Mid"
$ws.Cells.Item(5, 3).Value = "Mid"
$ws.Cells.Item(5, 4).Value = "High"
$ws.Cells.Item(6, 1).Value = "Apayao"
$ws.Cells.Item(6, 2).Value = "Low"
$ws.Cells.Item(6, 3).Value = "Low"
$ws.Cells.Item(6, 4).Value = "Low"
$ws.Cells.Item(7, 1).Value = "Aurora"
$ws.Cells.Item(7, 2).Value = "Mid"
$ws.Cells.Item(7, 3).Value = "Low"
$ws.Cells.Item(7, 4).Value = "Mid"
$ws.Cells.Item(8, 1).Value = "Basilan"
$ws.Cells.Item(8, 2).Value = "Low"
$ws.Cells.Item(8, 3).Value = "Low"
$ws.Cells.Item(8, 4).Value = "Low"
$ws.Cells.Item(9, 1).Value = "Bataan"
$ws.Cells.Item(9, 2).Value = "Low"
$ws.Cells.Item(9, 3).Value = "High"
$ws.Cells.Item(9, 4).Value = "Low"
$ws.Cells.Item(10, 1).Value = "Batanes"
$ws.Cells.Item(10, 2).Value = "Low"
$ws.Cells.Item(10, 3).Value = "Mid"
$ws.Cells.Item(10, 4).Value = "Low"
$ws.Cells.Item(11, 1).Value = "Batangas"
$ws.Cells.Item(11, 2).Value = "Low"
$ws.Cells.Item(11, 3).Value = "High"
$ws.Cells.Item(11, 4).Value = "Low"
$ws.Cells.Item(12, 1).Value = "Biliran"
$ws.Cells.Item(12, 2).Value = "High"
$ws.Cells.Item(12, 3).Value = "Mid"
$ws.Cells.Item(12, 4).Value = "High"
$ws.Cells.Item(13, 1).Value = "Bohol"
$ws.Cells.Item(13, 2).Value = "Mid"
$ws.Cells.Item(13, 3).Value = "Mid"
$ws.Cells.Item(13, 4).Value = "Mid"
$ws.Cells.Item(14, 1).Value = "Bukidnon"
$ws.Cells.Item(14, 2).Value = "Low"
$ws.Cells.Item(14, 3).Value = "Low"
$ws.Cells.Item(14, 4).Value = "Mid"
$ws.Cells.Item(15, 1).Value = "Bulacan"
$ws.Cells.Item(15, 2).Value = "Low"
$ws.Cells.Item(15, 3).Value = "High"
$ws.Cells.Item(15, 4).Value = "Low"
$ws.Cells.Item(16, 1).Value = "Cagayan"
$ws.Cells.Item(16, 2).Value = "High"
$ws.Cells.Item(16, 3).Value = "Mid"
$ws.Cells.Item(16, 4).Value = "High"
$ws.Cells.Item(17, 1).Value = "Camarines Norte"
$ws.Cells.Item(17, 2).Value = "High"
$ws.Cells.Item(17, 3).Value = "Mid"
$ws.Cells.Item(17, 4).Value = "High"
$ws.Cells.Item(18, 1).Value = "Camarines Sur"
$ws.Cells.Item(18, 2).Value = "High"
$ws.Cells.Item(18, 3).Value = "Mid"
$ws.Cells.Item(18, 4).Value = "High"
$ws.Cells.Item(19, 1).Value = "Camiguin"
$ws.Cells.Item(19, 2).Value = "Low"
$ws.Cells.Item(19, 3).Value = "Mid"
$ws.Cells.Item(19, 4).Value = "Mid"
$ws.Cells.Item(20, 1).Value = "Capiz"
$ws.Cells.Item(20, 2).Value = "High"
$ws.Cells.Item(20, 3).Value = "High"
$ws.Cells.Item(20, 4).Value = "High"
$ws.Cells.Item(21, 1).Value = "Catanduanes"
$ws.Cells.Item(21, 2).Value = "Mid"
$ws.Cells.Item(21, 3).Value = "Low"
$ws.Cells.Item(21, 4).Value = "Mid"
$ws.Cells.Item(22, 1).Value = "Cavite"
$ws.Cells.Item(22, 2).Value = "Mid"
$ws.Cells.Item(22, 3).Value = "High"
$ws.Cells.Item(22, 4).Value = "Low"
$ws.Cells.Item(23, 1).Value = "Cebu"
$ws.Cells.Item(23, 2).Value = "Mid"
$ws.Cells.Item(23, 3).Value = "High"
$ws.Cells.Item(23, 4).Value = "Mid"
$ws.Cells.Item(24, 1).Value = "Compostela Valley"
$ws.Cells.Item(24, 2).Value = "Low"
$ws.Cells.Item(24, 3).Value = "Mid"
$ws.Cells.Item(24, 4).Value = "Low"
$ws.Cells.Item(25, 1).Value = "Davao Oriental"
$ws.Cells.Item(25, 2).Value = "Low"
$ws.Cells.Item(25, 3).Value = "Low"
$ws.Cells.Item(25, 4).Value = "Mid"
$ws.Cells.Item(26, 1).Value = "Eastern Samar"
$ws.Cells.Item(26, 2).Value = "Mid"
$ws.Cells.Item(26, 3).Value = "Low"
$ws.Cells.Item(26, 4).Value = "High"
$ws.Cells.Item(27, 1).Value = "Guimaras"
$ws.Cells.Item(27, 2).Value = "High"
$ws.Cells.Item(27, 3).Value = "High"
$ws.Cells.Item(27, 4).Value = "High"
$ws.Cells.Item(28, 1).Value = "Ifugao"
$ws.Cells.Item(28, 2).Value = "Mid"
$ws.Cells.Item(28, 3).Value = "High"
$ws.Cells.Item(28, 4).Value = "Mid"
$ws.Cells.Item(29, 1).Value = "Ilocos Norte"
$ws.Cells.Item(29, 2).Value = "Mid"
$ws.Cells.Item(29, 3).Value = "High"
$ws.Cells.Item(29, 4).Value = "Low"
$ws.Cells.Item(30, 1).Value = "Ilocos Sur"
$ws.Cells.Item(30, 2).Value = "Mid"
$ws.Cells.Item(30, 3).Value = "High"
$ws.Cells.Item(30, 4).Value = "Low"
$ws.Cells.Item(31, 1).Value = "Isabela"
$ws.Cells.Item(31, 2).Value = "Mid"
$ws.Cells.Item(31, 3).Value = "Mid"
$ws.Cells.Item(31, 4).Value = "Mid"
$ws.Cells.Item(32, 1).Value = "Kalinga"
$ws.Cells.Item(32, 2).Value = "Mid"
$ws.Cells.Item(32, 3).Value = "Mid"
$ws.Cells.Item(32, 4).Value = "Mid"
$ws.Cells.Item(33, 1).Value = "La Union"
$ws.Cells.Item(33, 2).Value = "Low"
$ws.Cells.Item(33, 3).Value = "High"
$ws.Cells.Item(33, 4).Value = "Low"
$ws.Cells.Item(34, 1).Value = "Laguna"
$ws.Cells.Item(34, 2).Value = "Low"
$ws.Cells.Item(34, 3).Value = "High"
$ws.Cells.Item(34, 4).Value = "Low"
$ws.Cells.Item(35, 1).Value = "Leyte"
$ws.Cells.Item(35, 2).Value = "High"
$ws.Cells.Item(35, 3).Value = "Mid"
$ws.Cells.Item(35, 4).Value = "High"
$ws.Cells.Item(36, 1).Value = "Maguindanao"
$ws.Cells.Item(36, 2).Value = "High"
$ws.Cells.Item(36, 3).Value = "Low"
$ws.Cells.Item(36, 4).Value = "High"
$ws.Cells.Item(37, 1).Value = "Marinduque"
$ws.Cells.Item(37, 2).Value = "Mid"
$ws.Cells.Item(37, 3).Value = "Mid"
$ws.Cells.Item(37, 4).Value = "Mid"
$ws.Cells.Item(38, 1).Value = "Masbate"
$ws.Cells.Item(38, 2).Value = "Mid"
$ws.Cells.Item(38, 3).Value = "Low"
$ws.Cells.Item(38, 4).Value = "High"
$ws.Cells.Item(39, 1).Value = "Misamis Occidental"
$ws.Cells.Item(39, 2).Value = "Mid"
$ws.Cells.Item(39, 3).Value = "Low"
$ws.Cells.Item(39, 4).Value = "High"
$ws.Cells.Item(40, 1).Value = "Misamis Oriental"
$ws.Cells.Item(40, 2).Value = "High"
$ws.Cells.Item(40, 3).Value = "High"
$ws.Cells.Item(40, 4).Value = "Mid"
$ws.Cells.Item(41, 1).Value = "Negros Occidental"
$ws.Cells.Item(41, 2).Value = "High"
$ws.Cells.Item(41, 3).Value = "Mid"
$ws.Cells.Item(41, 4).Value = "High"
$ws.Cells.Item(42, 1).Value = "Negros Oriental"
$ws.Cells.Item(42, 2).Value = "Low"
$ws.Cells.Item(42, 3).Value = "Low"
$ws.Cells.Item(42, 4).Value = "Mid"
$ws.Cells.Item(43, 1).Value = "Northern Samar"
$ws.Cells.Item(43, 2).Value = "High"
$ws.Cells.Item(43, 3).Value = "Low"
$ws.Cells.Item(43, 4).Value = "High"
$ws.Cells.Item(44, 1).Value = "Nueva Ecija"
$ws.Cells.Item(44, 2).Value = "High"
$ws.Cells.Item(44, 3).Value = "High"
$ws.Cells.Item(44, 4).Value = "High"
$ws.Cells.Item(45, 1).Value = "Nueva Vizcaya"
$ws.Cells.Item(45, 2).Value = "High"
$ws.Cells.Item(45, 3).Value = "High"
$ws.Cells.Item(45, 4).Value = "High"
$ws.Cells.Item(46, 1).Value = "Occidental Mindoro"
$ws.Cells.Item(46, 2).Value = "Mid"
$ws.Cells.Item(46, 3).Value = "Mid"
$ws.Cells.Item(46, 4).Value = "Mid"
$ws.Cells.Item(47, 1).Value = "Oriental Mindoro"
$ws.Cells.Item(47, 2).Value = "Mid"
$ws.Cells.Item(47, 3).Value = "Mid"
$ws.Cells.Item(47, 4).Value = "Mid"
$ws.Cells.Item(48, 1).Value = "Palawan"
$ws.Cells.Item(48, 2).Value = "High"
$ws.Cells.Item(48, 3).Value = "Mid"
$ws.Cells.Item(48, 4).Value = "High"
$ws.Cells.Item(49, 1).Value = "Pampanga"
$ws.Cells.Item(49, 2).Value = "High"
$ws.Cells.Item(49, 3).Value = "High"
$ws.Cells.Item(49, 4).Value = "Mid"
$ws.Cells.Item(50, 1).Value = "Pangasinan"
$ws.Cells.Item(50, 2).Value = "High"
$ws.Cells.Item(50, 3).Value = "Mid"
$ws.Cells.Item(50, 4).Value = "High"
$ws.Cells.Item(51, 1).Value = "Quezon"
$ws.Cells.Item(51, 2).Value = "Mid"
$ws.Cells.Item(51, 3).Value = "Low"
$ws.Cells.Item(51, 4).Value = "Mid"
$ws.Cells.Item(52, 1).Value = "Rizal"
$ws.Cells.Item(52, 2).Value = "Low"
$ws.Cells.Item(52, 3).Value = "High"
$ws.Cells.Item(52, 4).Value = "Low"
$ws.Cells.Item(53, 1).Value = "Romblon"
$ws.Cells.Item(53, 2).Value = "Mid"
$ws.Cells.Item(53, 3).Value = "Low"
$ws.Cells.Item(53, 4).Value = "High"
$ws.Cells.Item(54, 1).Value = "Sarangani"
$ws.Cells.Item(54, 2).Value = "Low"
$ws.Cells.Item(54, 3).Value = "Low"
$ws.Cells.Item(54, 4).Value = "Mid"
$ws.Cells.Item(55, 1).Value = "Siquijor"
$ws.Cells.Item(55, 2).Value = "Mid"
$ws.Cells.Item(55, 3).Value = "High"
$ws.Cells.Item(55, 4).Value = "Mid"
$ws.Cells.Item(56, 1).Value = "Sorsogon"
$ws.Cells.Item(56, 2).Value = "Mid"
$ws.Cells.Item(56, 3).Value = "Low"
$ws.Cells.Item(56, 4).Value = "Mid"
$ws.Cells.Item(57, 1).Value = "South Cotabato"
$ws.Cells.Item(57, 2).Value = "Low"
$ws.Cells.Item(57, 3).Value = "Mid"
$ws.Cells.Item(57, 4).Value = "Low"
$ws.Cells.Item(58, 1).Value = "Southern Leyte"
$ws.Cells.Item(58, 2).Value = "Low"
$ws.Cells.Item(58, 3).Value = "Low"
$ws.Cells.Item(58, 4).Value = "Low"
$ws.Cells.Item(59, 1).Value = "Sultan Kudarat"
$ws.Cells.Item(59, 2).Value = "High"
$ws.Cells.Item(59, 3).Value = "Low"
$ws.Cells.Item(59, 4).Value = "High"
$ws.Cells.Item(60, 1).Value = "Tarlac"
$ws.Cells.Item(60, 2).Value = "High"
$ws.Cells.Item(60, 3).Value = "High"
$ws.Cells.Item(60, 4).Value = "Mid"
$ws.Cells.Item(61, 1).Value = "Zambales"
$ws.Cells.Item(61, 2).Value = "Low"
$ws.Cells.Item(61, 3).Value = "High"
$ws.Cells.Item(61, 4).Value = "Low"
$ws.Cells.Item(62, 1).Value = "Zamboanga Sibugay"
$ws.Cells.Item(62, 2).Value = "Low"
$ws.Cells.Item(62, 3).Value = "Low"
$ws.Cells.Item(62, 4).Value = "Low"
